$wb = $excel.ActiveWorkbook

# ALC row 137: Cutting Edge of Culinary Quality / Magnesia Whetstone
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6667635.5
$ws.Range("I137").Value = 991.7083
$ws.Range("J137").Value = 33334212
$ws.Range("K137").Value = 2975.1249
$ws.Range("L137").Value = 100002636
$ws.Range("M137").Value = -425.1248999999998
$ws.Range("N137").Value = -100007736

# ARM row 32: Ingot We Trust / Steel Ingot
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9433.796
$ws.Range("I32").Value = 9716.319
$ws.Range("K32").Value = 9716.319
$ws.Range("M32").Value = -9429.319

# ARM row 108: Time to Fry / Deepgold Rail Frypan
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H108").Value = 23182.5
$ws.Range("J108").Value = 23182.5
$ws.Range("L108").Value = 23182.5
$ws.Range("N108").Value = -30862.5

# BSM row 20: Smelt and Dealt / Iron Ingot
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2393.5715
$ws.Range("J20").Value = 2334.8333
$ws.Range("L20").Value = 2334.8333
$ws.Range("N20").Value = -2828.8333

# BSM row 86: Through Thick and Thin / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 23811256
$ws.Range("I86").Value = 1741.1428
$ws.Range("J86").Value = 71430290
$ws.Range("K86").Value = 1741.1428
$ws.Range("L86").Value = 71430290
$ws.Range("M86").Value = -618.1428000000001
$ws.Range("N86").Value = -71432536

# BSM row 89: Piercing Eyes Deserve Piercing Shafts (L) / Adamantite Nugget
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 23811256
$ws.Range("I89").Value = 1741.1428
$ws.Range("J89").Value = 71430290
$ws.Range("K89").Value = 8705.714
$ws.Range("L89").Value = 357151450
$ws.Range("M89").Value = -3089.714
$ws.Range("N89").Value = -357162682

# CRP row 80: The Long Armillae of the Law / Hallowed Chestnut Armillae
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H80").Value = 25000
$ws.Range("J80").Value = 25000
$ws.Range("L80").Value = 25000
$ws.Range("N80").Value = -27246

# CRP row 83: Wooden Ambitions (L) / Hallowed Chestnut Armillae
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H83").Value = 25000
$ws.Range("J83").Value = 25000
$ws.Range("L83").Value = 75000
$ws.Range("N83").Value = -86232

# CRP row 97: Wood That You Could / Larch Bracelets
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 21233.363
$ws.Range("J97").Value = 21233.363
$ws.Range("L97").Value = 21233.363
$ws.Range("N97").Value = -23215.363

# CRP row 99: O Pine / Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1658.238
$ws.Range("I99").Value = 1152.0714
$ws.Range("J99").Value = 2670.5715
$ws.Range("K99").Value = 1152.0714
$ws.Range("L99").Value = 2670.5715
$ws.Range("M99").Value = 345.9286
$ws.Range("N99").Value = -5666.5715

# CRP row 102: The Ear Is the Way to the Heart / Persimmon Earrings
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H102").Value = 10120.5
$ws.Range("J102").Value = 10120.5
$ws.Range("L102").Value = 10120.5
$ws.Range("N102").Value = -14988.5

# CRP row 104: Putting Your Line on the Neck / Zelkova Necklace
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H104").Value = 26353.334
$ws.Range("J104").Value = 26353.334
$ws.Range("L104").Value = 26353.334
$ws.Range("N104").Value = -31595.334

# CRP row 109: Playing the Market / White Oak Necklace
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H109").Value = 21813.77
$ws.Range("I109").Value = 20259
$ws.Range("J109").Value = 21943.334
$ws.Range("K109").Value = 20259
$ws.Range("L109").Value = 21943.334
$ws.Range("M109").Value = -19219
$ws.Range("N109").Value = -24023.334

# CRP row 126: A Better Conductor / Red Pine Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1658.238
$ws.Range("I126").Value = 1152.0714
$ws.Range("J126").Value = 2670.5715
$ws.Range("K126").Value = 3456.2142
$ws.Range("L126").Value = 8011.7145
$ws.Range("M126").Value = -986.2142000000003
$ws.Range("N126").Value = -12951.7145

# CRP row 132: Hull Lotta Damage / Ginseng Lumber
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 31254236
$ws.Range("I132").Value = 62503756
$ws.Range("J132").Value = 4717.5
$ws.Range("K132").Value = 187511268
$ws.Range("L132").Value = 14152.5
$ws.Range("M132").Value = -187508738
$ws.Range("N132").Value = -19212.5

# CRP row 140: Spear Pressure / Claro Walnut Spear
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 42920
$ws.Range("J140").Value = 42920
$ws.Range("L140").Value = 42920
$ws.Range("N140").Value = -53280

# CRP row 141: No Greater Treasure / Claro Walnut Necklace of Gathering
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 112861.086
$ws.Range("J141").Value = 112861.086
$ws.Range("L141").Value = 112861.086
$ws.Range("N141").Value = -123221.086

# CUL row 39: Bloody Good Tart, This / Blood Currant Tart
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 503.29413
$ws.Range("J39").Value = 503.29413
$ws.Range("L39").Value = 1509.88239
$ws.Range("N39").Value = -2097.88239

# CUL row 118: Teetotally / Masala Chai
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 946.4400000000001
$ws.Range("J118").Value = 961.15
$ws.Range("L118").Value = 2883.45
$ws.Range("N118").Value = -5369.45

# CUL row 140: Sweet, Sweet Bean Juice / Mesquite Juice
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3215.7144
$ws.Range("I140").Value = 2627.5
$ws.Range("J140").Value = 4000
$ws.Range("K140").Value = 7882.5
$ws.Range("L140").Value = 12000
$ws.Range("M140").Value = -2702.5
$ws.Range("N140").Value = -22360

# GSM row 15: The Tusk at Hand / Fang Earrings
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 16310.286
$ws.Range("J15").Value = 16310.286
$ws.Range("L15").Value = 16310.286
$ws.Range("N15").Value = -16886.286

# GSM row 70: Sky Is the Limit / Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 24352.3
$ws.Range("I70").Value = 37614.945
$ws.Range("J70").Value = 4458.3335
$ws.Range("K70").Value = 37614.945
$ws.Range("L70").Value = 4458.3335
$ws.Range("M70").Value = -37344.945
$ws.Range("N70").Value = -4998.3335

# GSM row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 24352.3
$ws.Range("I73").Value = 37614.945
$ws.Range("J73").Value = 4458.3335
$ws.Range("K73").Value = 37614.945
$ws.Range("L73").Value = 4458.3335
$ws.Range("M73").Value = -36678.945
$ws.Range("N73").Value = -6330.3335

# GSM row 81: The Grander Temple / Dragon Fang Earrings
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 16310.286
$ws.Range("J81").Value = 16310.286
$ws.Range("L81").Value = 16310.286
$ws.Range("N81").Value = -18306.286

# GSM row 84: Man with a Dragon Earring (L) / Dragon Fang Earrings
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 16310.286
$ws.Range("J84").Value = 16310.286
$ws.Range("L84").Value = 48930.858
$ws.Range("N84").Value = -58914.858

# LTW row 80: Don't Sweat the Small Fry / Dragonskin Wristbands
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H80").Value = 21996
$ws.Range("J80").Value = 21996
$ws.Range("L80").Value = 21996
$ws.Range("N80").Value = -24242

# LTW row 83: It's All in the Wrists (L) / Dragonskin Wristbands
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H83").Value = 21996
$ws.Range("J83").Value = 21996
$ws.Range("L83").Value = 65988
$ws.Range("N83").Value = -77220

# LTW row 132: Tenets of Tanning / Silver Lobo Leather
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 13166723
$ws.Range("I132").Value = 6116.353
$ws.Range("J132").Value = 23820548
$ws.Range("K132").Value = 18349.059
$ws.Range("L132").Value = 71461644
$ws.Range("M132").Value = -15819.059
$ws.Range("N132").Value = -71466704

# WVR row 102: Don't Sweat the Role / Serge Turban of Crafting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 33700
$ws.Range("J102").Value = 33700
$ws.Range("L102").Value = 33700
$ws.Range("N102").Value = -40190

# WVR row 106: Cap It Off / Serge Knit Cap
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H106").Value = 31666.666
$ws.Range("J106").Value = 31666.666
$ws.Range("L106").Value = 31666.666
$ws.Range("N106").Value = -34190.666

# WVR row 109: Turban in Training / Brightlinen Turban of Crafting
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 28450
$ws.Range("J109").Value = 28450
$ws.Range("L109").Value = 28450
$ws.Range("N109").Value = -31224

# WVR row 122: Heavy Armoire / Dark Hempen Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2731.1904
$ws.Range("I122").Value = 2637.4666
$ws.Range("J122").Value = 2965.5
$ws.Range("K122").Value = 7912.399800000001
$ws.Range("L122").Value = 8896.5
$ws.Range("M122").Value = -5462.399800000001
$ws.Range("N122").Value = -13796.5

# WVR row 132: Comfy Cabins / Snow Cotton Cloth
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4485.1113
$ws.Range("I132").Value = 3072.6
$ws.Range("J132").Value = 6250.75
$ws.Range("K132").Value = 9217.799999999999
$ws.Range("L132").Value = 18752.25
$ws.Range("M132").Value = -6687.799999999999
$ws.Range("N132").Value = -23812.25
